$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Feuil1" to "Taxonomy_update_2021"
$ws.Name = "Taxonomy_update_2021"

# The sheet had two blank, unused rows (17 and 18) between the main table
# (rows 1-16) and the trailing "Pseudoscada / timna" rows (formerly 19-21).
# Removing that gap shifts the trailing rows up to 17-19 and shrinks the
# used range from A1:F21 down to A1:F19.
$ws.Rows("17:18").Delete()

# Bold the header row (row 1) so the column titles stand out.
$ws.Range("A1:F1").Font.Bold = $true

# Leave the cursor where the author left it after editing.
$ws.Range("C16").Select()
